$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 995
$ws.Range("I80").Value = 995
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2985
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1987
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 995
$ws.Range("I83").Value = 995
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 8955
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -3963
$ws.Range("N83").ClearContents()
$ws.Range("H88").Value = 1521.3334
$ws.Range("I88").Value = 649
$ws.Range("J88").Value = 1630.375
$ws.Range("K88").Value = 649
$ws.Range("L88").Value = 1630.375
$ws.Range("M88").Value = -243
$ws.Range("N88").Value = -2442.375
$ws.Range("H91").Value = 1521.3334
$ws.Range("I91").Value = 649
$ws.Range("J91").Value = 1630.375
$ws.Range("K91").Value = 649
$ws.Range("L91").Value = 1630.375
$ws.Range("M91").Value = 755
$ws.Range("N91").Value = -4438.375
$ws.Range("H98").Value = 1679.8
$ws.Range("I98").Value = 1533.1111
$ws.Range("K98").Value = 1533.1111
$ws.Range("M98").Value = -35.11110000000008
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 1679.8
$ws.Range("I122").Value = 1533.1111
$ws.Range("K122").Value = 4599.3333
$ws.Range("M122").Value = -2149.3333
$ws.Range("N122").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1486.6
$ws.Range("I74").Value = 1486.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1486.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -612.5999999999999
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1486.6
$ws.Range("I77").Value = 1486.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7433
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3065
$ws.Range("N77").ClearContents()
$ws.Range("H122").Value = 5185.6
$ws.Range("I122").Value = 2778.4
$ws.Range("K122").Value = 8335.200000000001
$ws.Range("M122").Value = -5885.200000000001
$ws.Range("N122").ClearContents()
$ws.Range("H139").Value = 91665.336
$ws.Range("J139").Value = 91665.336
$ws.Range("L139").Value = 91665.336
$ws.Range("N139").Value = -101945.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 28686.625
$ws.Range("I81").Value = 50000
$ws.Range("J81").Value = 25641.857
$ws.Range("K81").Value = 50000
$ws.Range("L81").Value = 25641.857
$ws.Range("M81").Value = -48939
$ws.Range("N81").Value = -27763.857
$ws.Range("H84").Value = 28686.625
$ws.Range("I84").Value = 50000
$ws.Range("J84").Value = 25641.857
$ws.Range("K84").Value = 150000
$ws.Range("L84").Value = 76925.571
$ws.Range("M84").Value = -144696
$ws.Range("N84").Value = -87533.571
$ws.Range("H86").Value = 1922.1111
$ws.Range("I86").Value = 1703.9565
$ws.Range("J86").Value = 3176.5
$ws.Range("K86").Value = 1703.9565
$ws.Range("L86").Value = 3176.5
$ws.Range("M86").Value = -580.9565
$ws.Range("N86").Value = -5422.5
$ws.Range("H89").Value = 1922.1111
$ws.Range("I89").Value = 1703.9565
$ws.Range("J89").Value = 3176.5
$ws.Range("K89").Value = 8519.782499999999
$ws.Range("L89").Value = 15882.5
$ws.Range("M89").Value = -2903.782499999999
$ws.Range("N89").Value = -27114.5
$ws.Range("H94").Value = 551254
$ws.Range("I94").Value = 551254
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 551254
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -550803
$ws.Range("N94").ClearContents()
$ws.Range("H135").Value = 37475
$ws.Range("J135").Value = 37475
$ws.Range("L135").Value = 37475
$ws.Range("N135").Value = -47615

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1595
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H51").Value = 17972.5
$ws.Range("I51").Value = 17972.5
$ws.Range("K51").Value = 17972.5
$ws.Range("M51").Value = -17236.5
$ws.Range("H61").Value = 17972.5
$ws.Range("I61").Value = 17972.5
$ws.Range("K61").Value = 17972.5
$ws.Range("M61").Value = -17624.5
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H113").Value = 1595
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2530.4443
$ws.Range("I122").Value = 2221.75
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6665.25
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4215.25
$ws.Range("N122").Value = -19900

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1232.75
$ws.Range("I114").Value = 782.75
$ws.Range("J114").Value = 1682.75
$ws.Range("K114").Value = 2348.25
$ws.Range("L114").Value = 5048.25
$ws.Range("M114").Value = 905.75
$ws.Range("N114").Value = -11556.25
$ws.Range("H129").Value = 2547.5386
$ws.Range("J129").Value = 1999.875
$ws.Range("L129").Value = 5999.625
$ws.Range("N129").Value = -15999.625
$ws.Range("H131").Value = 940.6667
$ws.Range("I131").Value = 940.6667
$ws.Range("K131").Value = 2822.0001
$ws.Range("M131").Value = 2217.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5000000
$ws.Range("I11").Value = 5000000
$ws.Range("K11").Value = 5000000
$ws.Range("M11").Value = -4999861
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 3399.8
$ws.Range("I122").Value = 3249.75
$ws.Range("K122").Value = 9749.25
$ws.Range("M122").Value = -7299.25
$ws.Range("N122").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H55").Value = 2288.8572
$ws.Range("I55").Value = 2162
$ws.Range("K55").Value = 2162
$ws.Range("M55").Value = -1989
$ws.Range("N55").ClearContents()
$ws.Range("H68").Value = 8650
$ws.Range("I68").Value = 4600
$ws.Range("K68").Value = 4600
$ws.Range("M68").Value = -3851
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 8650
$ws.Range("I71").Value = 4600
$ws.Range("K71").Value = 23000
$ws.Range("M71").Value = -19256
$ws.Range("N71").ClearContents()
$ws.Range("H100").Value = 3124.6667
$ws.Range("J100").Value = 2895
$ws.Range("L100").Value = 2895
$ws.Range("N100").Value = -3977
$ws.Range("H132").Value = 3950
$ws.Range("I132").Value = 3950
$ws.Range("K132").Value = 11850
$ws.Range("M132").Value = -9320

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 7601
$ws.Range("I43").Value = 751.25
$ws.Range("J43").Value = 35000
$ws.Range("K43").Value = 751.25
$ws.Range("L43").Value = 35000
$ws.Range("M43").Value = -602.25
$ws.Range("N43").Value = -35298
$ws.Range("H81").Value = 2350
$ws.Range("I81").Value = 2350
$ws.Range("K81").Value = 4700
$ws.Range("M81").Value = -3639
$ws.Range("H84").Value = 2350
$ws.Range("I84").Value = 2350
$ws.Range("K84").Value = 23500
$ws.Range("M84").Value = -18196
$ws.Range("H122").Value = 1457.3334
$ws.Range("I122").Value = 1376.6842
$ws.Range("K122").Value = 4130.0526
$ws.Range("M122").Value = -1680.0526
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 4006.625
$ws.Range("I126").Value = 2008
$ws.Range("J126").Value = 10002.5
$ws.Range("K126").Value = 6024
$ws.Range("L126").Value = 30007.5
$ws.Range("M126").Value = -3554
$ws.Range("N126").Value = -34947.5
